# feat: add 2022-Q1 data
#
# 1) Insert a brand new "2022-Q1" worksheet (fund-holding detail for the
#    quarter) right before the existing "总计" (totals) sheet. It is
#    created by duplicating "2021-Q4" (same column layout/styles) and then
#    overwriting its data, which is simpler than building formatting from
#    scratch.
# 2) Insert a new leading row into "总计" summarising 2022-Q1, pushing the
#    previously existing quarters down by one row.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# 1) New "2022-Q1" sheet with the fund-holding detail table.
# ---------------------------------------------------------------------
$src = $wb.Worksheets.Item("2021-Q4")
$totalSheet = $wb.Worksheets.Item("总计")
$src.Copy($totalSheet, $null)

$q1 = $wb.Worksheets.Item("2021-Q4 (2)")
$q1.Name = "2022-Q1"

# The source sheet has 15 data rows (rows 2-16); only 4 are needed here.
$q1.Rows("6:16").Delete()

$rows = @(
    @("001044", "嘉实新消费股票", "8.92", "80.25", "5.07", "0.4522", 9),
    @("010551", "淳厚欣颐一年持有期混合", "3.41", "81.14", "3.59", "0.1224", 2),
    @("007811", "淳厚信泽灵活配置混合A", "4.32", "74.11", "2.32", "0.1002", 8),
    @("007812", "淳厚信泽灵活配置混合C", "1.15", "74.11", "2.32", "0.0267", 8)
)

for ($r = 0; $r -lt $rows.Length; $r++) {
    $row = $rows[$r]
    $excelRow = 2 + $r

    $q1.Cells.Item($excelRow, 1).Value = $r

    # B:G hold text values even though several look numeric (fund size,
    # positions, values) - force a text number format before writing so
    # leading zeros / exact decimal strings survive, then drop the format
    # again so the cells end up styleless like the source data.
    $textRange = $q1.Range("B" + $excelRow + ":G" + $excelRow)
    $textRange.NumberFormat = "@"
    for ($c = 0; $c -lt 6; $c++) {
        $q1.Cells.Item($excelRow, 2 + $c).Value = [string]$row[$c]
    }
    $textRange.ClearFormats()

    # H (仓位排名) is a genuine number.
    $q1.Cells.Item($excelRow, 8).Value = $row[6]
}

# ---------------------------------------------------------------------
# 2) Insert the new 2022-Q1 summary row at the top of "总计"'s data.
# ---------------------------------------------------------------------
$total = $wb.Worksheets.Item("总计")
$total.Rows(2).Insert()
$total.Rows(2).ClearFormats()

# Match the look of the other data rows: bordered/centered style on
# column A only (copied from the row right below, which still holds the
# old top row of data).
$total.Range("A3").Copy()
$total.Range("A2").PasteSpecial(-4122)

$total.Cells.Item(2, 1).Value = 0
$total.Cells.Item(2, 2).Value = "2022-Q1"
$total.Cells.Item(2, 3).Value = 4
$total.Cells.Item(2, 4).Value = 0.7

# Renumber the 0-based index column for the rows that got pushed down.
$lastRow = $total.UsedRange.Rows.Count
for ($r = 3; $r -le $lastRow; $r++) {
    $total.Cells.Item($r, 1).Value = $r - 2
}
